# Start a sliver of prob1
#
# 1. Merge the "Homework Set 3&4" / ":" runs into a single run and drop
#    the stray _GoBack bookmark that used to straddle them.
# 2. Mark every picture run as NoProofing (<w:noProof/>) - this is what
#    Word does the moment it lays pictures out / spell-checks around them.
# 3. Re-plant the _GoBack bookmark where the cursor was left last: right
#    after the two images in the first problem's paragraph.

$d = $word.ActiveDocument

# --- 1. Merge "Homework Set 3&4" + ":" into one run, eat the old bookmark ---
$d.Content.Find.Execute("Homework Set 3&4:", $false, $false, $false, $false, `
    $false, $true, 1, $false, "Homework Set 3&4:", 2) | Out-Null

# --- 2. Flag every inline picture's run as NoProofing ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $d.InlineShapes.Item($i).Range.NoProofing = 1
}

# The lone floating (anchored) picture isn't reached through InlineShapes;
# tag it too via its own range.
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $d.Shapes.Item($i).TextFrame.TextRange.NoProofing = 1
}

# --- 3. Move _GoBack to sit right after the second picture in that paragraph ---
$d.Bookmarks("_GoBack").Range.Select() | Out-Null
$target = $d.InlineShapes.Item(1).Range
$collapsed = $d.Range($target.End, $target.End)
$d.Bookmarks.Add("_GoBack", $collapsed) | Out-Null
